# Update the two-digit multiplication answers throughout the document.
# Each old expression is unique in the document, so a straightforward
# Find/Replace (non-wildcard, whole document) for each pair is safe.

$d = $word.ActiveDocument

$pairs = @(
    @("83×50=4150", "41×45=1845"),
    @("53×87=4611", "23×64=1472"),
    @("29×34=986",  "45×58=2610"),
    @("49×26=1274", "81×17=1377"),
    @("57×72=4104", "63×58=3654"),
    @("13×30=390",  "47×81=3807"),
    @("80×75=6000", "52×96=4992"),
    @("48×71=3408", "39×65=2535"),
    @("37×69=2553", "13×32=416"),
    @("43×75=3225", "68×76=5168"),
    @("16×33=528",  "34×59=2006"),
    @("39×15=585",  "42×52=2184"),
    @("29×75=2175", "16×94=1504"),
    @("98×86=8428", "37×79=2923"),
    @("69×16=1104", "13×98=1274"),
    @("36×67=2412", "25×68=1700"),
    @("32×91=2912", "71×54=3834"),
    @("63×34=2142", "80×81=6480"),
    @("23×70=1610", "27×59=1593"),
    @("63×20=1260", "58×70=4060"),
    @("80×83=6640", "35×14=490"),
    @("79×92=7268", "36×37=1332"),
    @("96×48=4608", "69×39=2691"),
    @("99×39=3861", "96×69=6624"),
    @("47×34=1598", "63×56=3528")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
